$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 62.41592166666667
$ws.Cells.Item(2, 8).Value = 187.247765
$ws.Cells.Item(2, 9).Value = 0.1654944774607044
$ws.Cells.Item(2, 10).Value = 0.1654944774607044
$ws.Cells.Item(2, 13).Value = 0.6327629999999999
$ws.Cells.Item(2, 14).Value = 1.898289
$ws.Cells.Item(2, 15).Value = 0.1382544270550543
$ws.Cells.Item(2, 16).Value = 0.1382544270550544
$ws.Cells.Item(2, 17).Value = 39.49448584156499
$ws.Cells.Item(2, 18).Value = 355.450372574085
$ws.Cells.Item(2, 19).Value = 0.0228803441621053
$ws.Cells.Item(2, 20).Value = 0.0228803441621053
$ws.Cells.Item(3, 7).Value = 62.41592166666667
$ws.Cells.Item(3, 8).Value = 187.247765
$ws.Cells.Item(3, 9).Value = 0.1654944774607044
$ws.Cells.Item(3, 10).Value = 0.1654944774607044
$ws.Cells.Item(3, 15).Value = 0.4765301499162115
$ws.Cells.Item(3, 16).Value = 0.4765301499162115
$ws.Cells.Item(3, 17).Value = 136.1281056949461
$ws.Cells.Item(3, 18).Value = 1225.152951254515
$ws.Cells.Item(3, 19).Value = 0.07886310815465457
$ws.Cells.Item(3, 20).Value = 0.07886310815465458
$ws.Cells.Item(4, 7).Value = 62.41592166666667
$ws.Cells.Item(4, 8).Value = 187.247765
$ws.Cells.Item(4, 9).Value = 0.1654944774607044
$ws.Cells.Item(4, 10).Value = 0.1654944774607044
$ws.Cells.Item(4, 13).Value = 1.444396333333334
$ws.Cells.Item(4, 14).Value = 4.333189000000001
$ws.Cells.Item(4, 15).Value = 0.3155908096798033
$ws.Cells.Item(4, 16).Value = 0.3155908096798033
$ws.Cells.Item(4, 17).Value = 90.15332839695391
$ws.Cells.Item(4, 18).Value = 811.3799555725852
$ws.Cells.Item(4, 19).Value = 0.05222853613935968
$ws.Cells.Item(4, 20).Value = 0.05222853613935968
$ws.Cells.Item(5, 7).Value = 62.41592166666667
$ws.Cells.Item(5, 8).Value = 187.247765
$ws.Cells.Item(5, 9).Value = 0.1654944774607044
$ws.Cells.Item(5, 10).Value = 0.1654944774607044
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.3186579999999999
$ws.Cells.Item(5, 14).Value = 0.9559739999999999
$ws.Cells.Item(5, 15).Value = 0.06962461334893082
$ws.Cells.Item(5, 16).Value = 0.06962461334893082
$ws.Cells.Item(5, 17).Value = 19.88933276645666
$ws.Cells.Item(5, 18).Value = 179.00399489811
$ws.Cells.Item(5, 19).Value = 0.01152248900458489
$ws.Cells.Item(5, 20).Value = 0.01152248900458489
$ws.Cells.Item(6, 9).Value = 0.4369365253446571
$ws.Cells.Item(6, 10).Value = 0.436936525344657
$ws.Cells.Item(6, 13).Value = 0.6327629999999999
$ws.Cells.Item(6, 14).Value = 1.898289
$ws.Cells.Item(6, 15).Value = 0.1382544270550543
$ws.Cells.Item(6, 16).Value = 0.1382544270550544
$ws.Cells.Item(6, 17).Value = 104.272865648853
$ws.Cells.Item(6, 18).Value = 938.4557908396769
$ws.Cells.Item(6, 19).Value = 0.06040840897095179
$ws.Cells.Item(6, 20).Value = 0.0604084089709518
$ws.Cells.Item(7, 9).Value = 0.4369365253446571
$ws.Cells.Item(7, 10).Value = 0.436936525344657
$ws.Cells.Item(7, 15).Value = 0.4765301499162115
$ws.Cells.Item(7, 16).Value = 0.4765301499162115
$ws.Cells.Item(7, 19).Value = 0.208213427926358
$ws.Cells.Item(7, 20).Value = 0.208213427926358
$ws.Cells.Item(8, 9).Value = 0.4369365253446571
$ws.Cells.Item(8, 10).Value = 0.436936525344657
$ws.Cells.Item(8, 13).Value = 1.444396333333334
$ws.Cells.Item(8, 14).Value = 4.333189000000001
$ws.Cells.Item(8, 15).Value = 0.3155908096798033
$ws.Cells.Item(8, 16).Value = 0.3155908096798033
$ws.Cells.Item(8, 17).Value = 238.0217313739309
$ws.Cells.Item(8, 18).Value = 2142.195582365378
$ws.Cells.Item(8, 19).Value = 0.1378931518122002
$ws.Cells.Item(8, 20).Value = 0.1378931518122002
$ws.Cells.Item(9, 9).Value = 0.4369365253446571
$ws.Cells.Item(9, 10).Value = 0.436936525344657
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.3186579999999999
$ws.Cells.Item(9, 14).Value = 0.9559739999999999
$ws.Cells.Item(9, 15).Value = 0.06962461334893082
$ws.Cells.Item(9, 16).Value = 0.06962461334893082
$ws.Cells.Item(9, 17).Value = 52.51157672293133
$ws.Cells.Item(9, 18).Value = 472.6041905063819
$ws.Cells.Item(9, 19).Value = 0.03042153663514706
$ws.Cells.Item(9, 20).Value = 0.03042153663514706
$ws.Cells.Item(10, 7).Value = 57.486235
$ws.Cells.Item(10, 8).Value = 172.458705
$ws.Cells.Item(10, 9).Value = 0.1524235190071549
$ws.Cells.Item(10, 10).Value = 0.1524235190071549
$ws.Cells.Item(10, 13).Value = 0.6327629999999999
$ws.Cells.Item(10, 14).Value = 1.898289
$ws.Cells.Item(10, 15).Value = 0.1382544270550543
$ws.Cells.Item(10, 16).Value = 0.1382544270550544
$ws.Cells.Item(10, 17).Value = 36.37516251730499
$ws.Cells.Item(10, 18).Value = 327.3764626557449
$ws.Cells.Item(10, 19).Value = 0.02107322629004938
$ws.Cells.Item(10, 20).Value = 0.02107322629004939
$ws.Cells.Item(11, 7).Value = 57.486235
$ws.Cells.Item(11, 8).Value = 172.458705
$ws.Cells.Item(11, 9).Value = 0.1524235190071549
$ws.Cells.Item(11, 10).Value = 0.1524235190071549
$ws.Cells.Item(11, 15).Value = 0.4765301499162115
$ws.Cells.Item(11, 16).Value = 0.4765301499162115
$ws.Cells.Item(11, 17).Value = 125.3765395931617
$ws.Cells.Item(11, 18).Value = 1128.388856338455
$ws.Cells.Item(11, 19).Value = 0.07263440236323603
$ws.Cells.Item(11, 20).Value = 0.07263440236323604
$ws.Cells.Item(12, 7).Value = 57.486235
$ws.Cells.Item(12, 8).Value = 172.458705
$ws.Cells.Item(12, 9).Value = 0.1524235190071549
$ws.Cells.Item(12, 10).Value = 0.1524235190071549
$ws.Cells.Item(12, 13).Value = 1.444396333333334
$ws.Cells.Item(12, 14).Value = 4.333189000000001
$ws.Cells.Item(12, 15).Value = 0.3155908096798033
$ws.Cells.Item(12, 16).Value = 0.3155908096798033
$ws.Cells.Item(12, 17).Value = 83.03290705113835
$ws.Cells.Item(12, 18).Value = 747.2961634602452
$ws.Cells.Item(12, 19).Value = 0.0481034617777129
$ws.Cells.Item(12, 20).Value = 0.0481034617777129
$ws.Cells.Item(13, 7).Value = 57.486235
$ws.Cells.Item(13, 8).Value = 172.458705
$ws.Cells.Item(13, 9).Value = 0.1524235190071549
$ws.Cells.Item(13, 10).Value = 0.1524235190071549
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.3186579999999999
$ws.Cells.Item(13, 14).Value = 0.9559739999999999
$ws.Cells.Item(13, 15).Value = 0.06962461334893082
$ws.Cells.Item(13, 16).Value = 0.06962461334893082
$ws.Cells.Item(13, 17).Value = 18.31844867263
$ws.Cells.Item(13, 18).Value = 164.86603805367
$ws.Cells.Item(13, 19).Value = 0.01061242857615657
$ws.Cells.Item(13, 20).Value = 0.01061242857615657
$ws.Cells.Item(14, 7).Value = 92.45614233333333
$ws.Cells.Item(14, 8).Value = 277.368427
$ws.Cells.Item(14, 9).Value = 0.2451454781874835
$ws.Cells.Item(14, 10).Value = 0.2451454781874835
$ws.Cells.Item(14, 13).Value = 0.6327629999999999
$ws.Cells.Item(14, 14).Value = 1.898289
$ws.Cells.Item(14, 15).Value = 0.1382544270550543
$ws.Cells.Item(14, 16).Value = 0.1382544270550544
$ws.Cells.Item(14, 17).Value = 58.50282599126699
$ws.Cells.Item(14, 18).Value = 526.5254339214029
$ws.Cells.Item(14, 19).Value = 0.03389244763194785
$ws.Cells.Item(14, 20).Value = 0.03389244763194786
$ws.Cells.Item(15, 7).Value = 92.45614233333333
$ws.Cells.Item(15, 8).Value = 277.368427
$ws.Cells.Item(15, 9).Value = 0.2451454781874835
$ws.Cells.Item(15, 10).Value = 0.2451454781874835
$ws.Cells.Item(15, 15).Value = 0.4765301499162115
$ws.Cells.Item(15, 16).Value = 0.4765301499162115
$ws.Cells.Item(15, 17).Value = 201.6453363120085
$ws.Cells.Item(15, 18).Value = 1814.808026808077
$ws.Cells.Item(15, 19).Value = 0.1168192114719629
$ws.Cells.Item(15, 20).Value = 0.1168192114719629
$ws.Cells.Item(16, 7).Value = 92.45614233333333
$ws.Cells.Item(16, 8).Value = 277.368427
$ws.Cells.Item(16, 9).Value = 0.2451454781874835
$ws.Cells.Item(16, 10).Value = 0.2451454781874835
$ws.Cells.Item(16, 13).Value = 1.444396333333334
$ws.Cells.Item(16, 14).Value = 4.333189000000001
$ws.Cells.Item(16, 15).Value = 0.3155908096798033
$ws.Cells.Item(16, 16).Value = 0.3155908096798033
$ws.Cells.Item(16, 17).Value = 133.5433129804115
$ws.Cells.Item(16, 18).Value = 1201.889816823703
$ws.Cells.Item(16, 19).Value = 0.07736565995053049
$ws.Cells.Item(16, 20).Value = 0.07736565995053049
$ws.Cells.Item(17, 7).Value = 92.45614233333333
$ws.Cells.Item(17, 8).Value = 277.368427
$ws.Cells.Item(17, 9).Value = 0.2451454781874835
$ws.Cells.Item(17, 10).Value = 0.2451454781874835
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.3186579999999999
$ws.Cells.Item(17, 14).Value = 0.9559739999999999
$ws.Cells.Item(17, 15).Value = 0.06962461334893082
$ws.Cells.Item(17, 16).Value = 0.06962461334893082
$ws.Cells.Item(17, 17).Value = 29.46188940365533
$ws.Cells.Item(17, 18).Value = 265.157004632898
$ws.Cells.Item(17, 19).Value = 0.0170681591330423
$ws.Cells.Item(17, 20).Value = 0.0170681591330423
